$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 58

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-19"
$ws.Cells.Item($row, 1).ClearFormats()
$ws.Cells.Item($row, 2).Value = "15:57:55"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "25"
$ws.Cells.Item($row, 4).ClearFormats()
$ws.Cells.Item($row, 5).Value = 121956
$ws.Cells.Item($row, 6).Value = 133821
$ws.Cells.Item($row, 7).Value = 162264
$ws.Cells.Item($row, 8).Value = 133093
$ws.Cells.Item($row, 9).Value = 177306
$ws.Cells.Item($row, 10).Value = 114598
$ws.Cells.Item($row, 11).Value = 201468
$ws.Cells.Item($row, 12).Value = 225031
$ws.Cells.Item($row, 13).Value = 175468
$ws.Cells.Item($row, 14).Value = 103766
$ws.Cells.Item($row, 15).Value = 39145
$ws.Cells.Item($row, 16).Value = 33944
$ws.Cells.Item($row, 17).Value = 51831
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36346
$ws.Cells.Item($row, 20).Value = -1
